# Corrected incremental heat rate
# Adds a new "Steam_Plant" unit row to Table1 on the "Units" sheet and
# updates the active sheet/selection to match the saved workbook view.

$wb = $excel.ActiveWorkbook

# --- Add a new row to Table1 ("Units" sheet) for the Steam_Plant unit ---
$ws = $wb.Worksheets.Item("Units")
$tbl = $ws.ListObjects.Item("Table1")
$newRow = $tbl.ListRows.Add()

$ws.Range("A7").Value = "Steam_Plant"   # Unit
$ws.Range("B7").Value = "Water"         # Input1
$ws.Range("C7").Value = "Waste_Heat"    # Input2
$ws.Range("D7").Value = "Steam"         # Output1
$ws.Range("J7").Value = 200             # Cap_Output1_existing
$ws.Range("W7").Value = 0.2             # Relation_In_In

# Calculated "Error messages:" column formula for the new row
$ws.Range("AJ7").Formula = '=IF( Table1[[#This Row],[minimum_op_point]]="", "", IF( COUNTA(Table1[[#This Row],[Cap_Input1_existing]], Table1[[#This Row],[Cap_Input2_existing]], Table1[[#This Row],[Cap_Output1_existing]], Table1[[#This Row],[Cap_Output2_existing]]) = 1, "", IF( COUNTA(Table1[[#This Row],[Cap_Input1_existing]], Table1[[#This Row],[Cap_Input2_existing]], Table1[[#This Row],[Cap_Output1_existing]], Table1[[#This Row],[Cap_Output2_existing]]) = 0, "Capacity missing", "Too many capacities" ) ) )'

# Match the red-font error-column style used elsewhere in the table
$ws.Range("AJ1").Copy()
$ws.Range("AJ7").PasteSpecial(-4122)

# --- Update the active sheet / selection to match saved workbook state ---
$ws.Activate() | Out-Null
$ws.Range("AG1").Select() | Out-Null
